$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Análisis" (Analysis) section values in column A (rows 20-24)
$ws.Range("A20").Value = 2
$ws.Range("A21").Value = 1
$ws.Range("A22").Value = 1
$ws.Range("A23").Value = 1
$ws.Range("A24").Value = 1

# Update the selected cell/active cell to D25
$ws.Range("D25").Select()
